# Flujo originación con compra cartera hasta analisis
# Update the "OriginacionDigiCredito" worksheet (3rd sheet) to add a
# cartera/saneamiento pair of columns (Cartera1 / Saneamiento2) with their
# sample data, and change the sample "Monto" value in F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OriginacionDigiCredito")
$ws.Activate()

# F2: Monto sample value changes from "20000000" to "7500000"
$ws.Range("F2").Value = '"7500000"'

# New header columns AY1/AZ1
$ws.Range("AY1").Value = "Cartera1"
$ws.Range("AZ1").Value = "Saneamiento2"

# New sample data columns AY2/AZ2
$ws.Range("AY2").Value = '"200000"'
$ws.Range("AZ2").Value = '"250000"'

# Move the active selection to the last new header cell, matching the
# author's view state when the change was made.
$ws.Range("AX1").Select()
